$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of trade data (row 6)
$row = 6

# A6: date/time serial value, formatted like A3:A5 (numFmtId 22 date-time) -> copy style from A5
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A$row").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A$row").Value = 42647.680671296293

# B6: boolean TRUE
$ws.Range("B$row").Value = $true

# C6 - F6: numeric values
$ws.Range("C$row").Value = 9941.89
$ws.Range("D$row").Value = 9766.58
$ws.Range("E$row").Value = 18.12
$ws.Range("F$row").Value = 18.77

# G6: boolean FALSE (column already styled via column default)
$ws.Range("G$row").Value = $false

# H6: numeric value
$ws.Range("H$row").Value = 3.59

# I6: boolean FALSE
$ws.Range("I$row").Value = $false

$excel.CutCopyMode = $false
